# Updates the cryptos price/volume table with freshly scraped figures.
# All D/E columns in the sheet are stored as text (prices use dotted
# grouping like "25.696.39" and volumes keep padding like "  -0.35%  "),
# so every write below targets Range.Value with a string. A handful of
# the new Price values look like plain numbers (e.g. "214.44") and Excel
# would silently reinterpret those as numeric cells; a leading "'" forces
# the same text storage the workbook already uses, exactly like typing an
# apostrophe-prefixed entry in the UI (the marker itself is not stored).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.696.39'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.628.76'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''214.44'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '''0.254'
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('D9').Value = '''0.0631'
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').Value = '''19.43'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('E11').Value = '  +1.34%  '
$ws.Range('D12').Value = '''4.25'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '1.856.14'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').Value = '1.626.26'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').Value = '''0.553'
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = '0.0₃0759'
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('D17').Value = '''62.75'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').Value = '25.710.99'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = '''4.43'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Value = '''191.76'
$ws.Range('E21').Value = '  -1.37%  '
$ws.Range('D22').Value = '''9.91'
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').Value = '''6.24'
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  +3.28%  '
$ws.Range('D26').Value = '''142.20'
$ws.Range('E26').Value = '  +2.08%  '
$ws.Range('E27').Value = '  +1.52%  '
$ws.Range('D28').Value = '''6.84'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('D29').Value = '''15.44'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('D30').Value = '''1.23'
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('E32').Value = '  -0.71%  '
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('D35').Value = '''2.39'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('D37').Value = '1.141.56'
$ws.Range('E37').Value = '  +3.04%  '
$ws.Range('E38').Value = '  -2.57%  '
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('B42').Value = 'mCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D42').Value = '''2.54'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '''5.53'
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('D44').Value = '''100.31'
$ws.Range('E44').Value = '  +1.06%  '
$ws.Range('D45').Value = '''0.802'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').Value = '1.764.87'
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('D47').Value = '0.0₆0110'
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('D48').Value = '''55.17'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('E49').Value = '  +0.79%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').Value = '''1.44'
$ws.Range('E51').Value = '  +4.51%  '
